$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.464.56'
$ws.Range("E2").Value = '  +1.79%  '
$ws.Range("D3").Value = '2.164.01'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '228.70'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.00%  '
$ws.Range("E6").Value = '  +1.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '63.78'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.20%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  +1.64%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0854'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.31%  '
$ws.Range("E11").Value = '  +0.41%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.08'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.02%  '
$ws.Range("D13").Value = '2.486.18'
$ws.Range("E13").Value = '  +2.90%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.11'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.08%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.814'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.67%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.53'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.28%  '
$ws.Range("D17").Value = '2.165.74'
$ws.Range("E17").Value = '  +2.89%  '
$ws.Range("D18").Value = '39.453.93'
$ws.Range("E18").Value = '  +1.70%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.19'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.46%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.90'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = '0.0₃0850'
$ws.Range("E21").Value = '  +0.99%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '229.56'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.61%  '
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.35'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("E25").Value = '  +1.82%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '172.26'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.45%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.52'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.43%  '
$ws.Range("E28").Value = '  +1.63%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.88'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.66%  '
$ws.Range("E30").Value = '  -0.33%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.66'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.95%  '
$ws.Range("E32").Value = '  +1.09%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.63'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.76%  '
$ws.Range("E34").Value = '  +1.78%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.72'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.99%  '
$ws.Range("E36").Value = '  +0.46%  '
$ws.Range("E37").Value = '  +1.05%  '
$ws.Range("E38").Value = '  -0.05%  '
$ws.Range("E39").Value = '  +0.23%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '103.19'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.30%  '
$ws.Range("E41").Value = '  +0.93%  '
$ws.Range("D43").Value = '1.524.41'
$ws.Range("E43").Value = '  -0.83%  '
$ws.Range("E44").Value = '  +4.39%  '
$ws.Range("E45").Value = '  +5.81%  '
$ws.Range("E46").Value = '  +2.06%  '
$ws.Range("E47").Value = '  +0.85%  '
$ws.Range("E48").Value = '  +4.11%  '
$ws.Range("E49").Value = '  -1.93%  '
$ws.Range("D50").Value = '2.368.91'
$ws.Range("E50").Value = '  +2.90%  '
$ws.Range("E51").Value = '  -0.69%  '
